$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates between rows 4-5 and rows 6-7 in column D.
# Rows 4 and 5 were 44574 (2022-01-13) -> should become 44559 (2021-12-29)
# Rows 6 and 7 were 44559 (2021-12-29) -> should become 44574 (2022-01-13)
$ws.Range("D4").Value = 44559
$ws.Range("D5").Value = 44559
$ws.Range("D6").Value = 44574
$ws.Range("D7").Value = 44574
